$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.316.12"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.819.88"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.44%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.48"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4355"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.95%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3690"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.40%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.00"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07717"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.141"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.59%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.24"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.001"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.343"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.579"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +5.83%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.837.29"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.60%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.54"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +7.46%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06525"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +8.35%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.02%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.291"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.359.14"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.70"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.011"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -14.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.35"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.33%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.038.69"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.305"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.89%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.14"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.84%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.218"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.012"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.22%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09221"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.11%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.550"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.93%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.00"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02365"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.57%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2186"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.70%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.221"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6618"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06215"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.188"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.58%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.52%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.441"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.72%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9996"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.08"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6147"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.760"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.69%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.60"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.99%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.035"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.94%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.161"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.37%  "

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07024"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.03%  "
